# Daily attendance processing - 2025-10-27 15:44:16
#
# Reorders the comma-separated list of editors/users in the "Recorded By"
# column (G) for specific rows on the "Session Analysis Results" sheet.
#
# Rule observed in the target diff:
#   - If the first token is the exact (case-sensitive) word "System",
#     rotate the list so that token moves to the end
#     (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System").
#   - If the first token is lower-case "system" (kept in place), the
#     exact token "System" elsewhere in the list is moved to the end
#     (e.g. "system, System, backup@backdoor.com" ->
#           "system, backup@backdoor.com, System").
#   - Otherwise (no exact "System" token present), simply rotate the
#     first token to the end
#     (e.g. "admin@admin.com, dnasr281@gmail.com" ->
#           "dnasr281@gmail.com, admin@admin.com").

function Transform-RecordedBy($s) {
    $rawParts = $s.Split(",")
    $parts = New-Object System.Collections.ArrayList
    foreach ($p in $rawParts) {
        $parts.Add($p.Trim()) | Out-Null
    }

    if ($parts.Count -lt 2) {
        return $s
    }

    if ($parts[0].Equals("system")) {
        # Lower-case "system" stays first; move the exact token "System"
        # (if present among the remaining tokens) to the very end.
        $rest = New-Object System.Collections.ArrayList
        for ($i = 1; $i -lt $parts.Count; $i++) {
            $rest.Add($parts[$i]) | Out-Null
        }

        $foundIndex = -1
        for ($i = 0; $i -lt $rest.Count; $i++) {
            if ($rest[$i].Equals("System")) {
                $foundIndex = $i
            }
        }

        $newParts = New-Object System.Collections.ArrayList
        $newParts.Add($parts[0]) | Out-Null
        if ($foundIndex -ge 0) {
            for ($i = 0; $i -lt $rest.Count; $i++) {
                if ($i -ne $foundIndex) {
                    $newParts.Add($rest[$i]) | Out-Null
                }
            }
            $newParts.Add("System") | Out-Null
        } else {
            for ($i = 0; $i -lt $rest.Count; $i++) {
                $newParts.Add($rest[$i]) | Out-Null
            }
        }
        return [string]::Join(", ", $newParts)
    } else {
        # Rotate: move the first token to the end of the list.
        $newParts = New-Object System.Collections.ArrayList
        for ($i = 1; $i -lt $parts.Count; $i++) {
            $newParts.Add($parts[$i]) | Out-Null
        }
        $newParts.Add($parts[0]) | Out-Null
        return [string]::Join(", ", $newParts)
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Exact rows touched by the commit (column G = "Recorded By").
$targetRows = @(
    2, 3, 5, 6, 7, 8, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24,
    29, 30, 32, 33, 34, 35, 37, 38, 39, 40, 41, 42, 44, 45, 46, 47, 48, 49, 51,
    56, 57, 59, 60, 61, 62, 64, 65, 66, 67, 68, 69, 71, 72, 73, 74, 75, 76, 78,
    83, 84, 85, 86, 87, 88, 89, 90, 93, 95, 96, 97, 99, 102,
    109, 110, 111, 112, 113, 114, 115, 116, 119, 121, 122, 123, 125, 128,
    135, 136, 137, 138, 139, 140, 141, 142, 145, 147, 148, 149, 151, 154
)

foreach ($row in $targetRows) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2
    $updated = Transform-RecordedBy $current
    $cell.Value = $updated
}
